# Fix js calendar. Fix Bożena and 'Name, Surname' in resources
#
# 1) "Iwaniec Joanna" -> "Joanna Iwaniec" (Name, Surname fix) everywhere it is used
# 2) "Giermek Bozena" (ASCII typo duplicate) -> "Giermek Bożena" (correct existing
#    string with the Polish diacritic) everywhere it is used, so the bad duplicate
#    shared string disappears from the workbook entirely
# 3) Selection / active cell bookkeeping left behind by the edit session

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Fix "Iwaniec Joanna" -> "Joanna Iwaniec" -----------------------------
$joannaCells = @("E2","E4","E7","E21","E23","E24","E26","E41","E44","E48","E58","E59","E61","E64","E68","E71","E76")
foreach ($addr in $joannaCells) {
    $ws1.Range($addr).Value2 = "Joanna Iwaniec"
}

# --- Fix "Giermek Bozena" -> "Giermek Bożena" (drop duplicate string) -----
$bozenaCells = @("E9","E25","E63","E66")
foreach ($addr in $bozenaCells) {
    $ws1.Range($addr).Value2 = "Giermek Bożena"
}

# --- Restore selection state on each sheet --------------------------------
$ws1.Activate()
$ws1.Range("E76").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A1").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("A1").Select()

$ws1.Activate()
